$wb = $excel.ActiveWorkbook

# ===== ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 4319.8335
$ws.Range("I7").Value = 3833.3333
$ws.Range("K7").Value = 3833.3333
$ws.Range("M7").Value = -3721.3333
$ws.Range("H14").Value = 4319.8335
$ws.Range("I14").Value = 3833.3333
$ws.Range("K14").Value = 3833.3333
$ws.Range("M14").Value = -3642.3333
$ws.Range("H39").Value = 173.21053
$ws.Range("J39").Value = 498.5
$ws.Range("L39").Value = 1495.5
$ws.Range("N39").Value = -2087.5
$ws.Range("H43").Value = 5375.857
$ws.Range("I43").Value = 3666.3333
$ws.Range("J43").Value = 5842.091
$ws.Range("K43").Value = 3666.3333
$ws.Range("L43").Value = 5842.091
$ws.Range("M43").Value = -3597.3333
$ws.Range("N43").Value = -5980.091
$ws.Range("H62").Value = 68217.69500000001
$ws.Range("I62").Value = 103353.75
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 103353.75
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = -102729.75
$ws.Range("N62").Value = -13248
$ws.Range("H64").Value = 8666.619000000001
$ws.Range("I64").Value = 4000
$ws.Range("K64").Value = 4000
$ws.Range("M64").Value = -3752
$ws.Range("H65").Value = 68217.69500000001
$ws.Range("I65").Value = 103353.75
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 516768.75
$ws.Range("L65").Value = 60000
$ws.Range("M65").Value = -513648.75
$ws.Range("N65").Value = -66240
$ws.Range("H67").Value = 8666.619000000001
$ws.Range("I67").Value = 4000
$ws.Range("K67").Value = 4000
$ws.Range("M67").Value = -3142
$ws.Range("H98").Value = 5737.4194
$ws.Range("I98").Value = 7791.864
$ws.Range("J98").Value = 715.44446
$ws.Range("K98").Value = 7791.864
$ws.Range("L98").Value = 715.44446
$ws.Range("M98").Value = -6293.864
$ws.Range("N98").Value = -3711.44446
$ws.Range("H112").Value = 1964.7693
$ws.Range("I112").Value = 1249.4445
$ws.Range("K112").Value = 3748.3335
$ws.Range("M112").Value = -2640.3335
$ws.Range("H122").Value = 5737.4194
$ws.Range("I122").Value = 7791.864
$ws.Range("J122").Value = 715.44446
$ws.Range("K122").Value = 23375.592
$ws.Range("L122").Value = 2146.33338
$ws.Range("M122").Value = -20925.592
$ws.Range("N122").Value = -7046.33338
$ws.Range("H137").Value = 2072.4211
$ws.Range("J137").Value = 2396.8
$ws.Range("L137").Value = 7190.400000000001
$ws.Range("N137").Value = -12290.4
$ws.Range("H138").Value = 1718.091
$ws.Range("I138").Value = 1097.1892
$ws.Range("K138").Value = 3291.5676
$ws.Range("M138").Value = 1848.4324
$ws.Range("H141").Value = 4903
$ws.Range("I141").Value = 5490.8335
$ws.Range("K141").Value = 16472.5005
$ws.Range("M141").Value = -11292.5005

# ===== ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2474.0193
$ws.Range("I32").Value = 2522.2744
$ws.Range("K32").Value = 2522.2744
$ws.Range("M32").Value = -2235.2744
$ws.Range("H45").Value = 1190
$ws.Range("I45").Value = 1190
$ws.Range("K45").Value = 1190
$ws.Range("M45").Value = -813
$ws.Range("H61").Value = 4801.1875
$ws.Range("I61").Value = 4801.1875
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4801.1875
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -4589.1875
$ws.Range("H110").Value = 1016.6667
$ws.Range("I110").Value = 1016.6667
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1016.6667
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = 1028.3333
$ws.Range("H122").Value = 10103175
$ws.Range("I122").Value = 12347692
$ws.Range("K122").Value = 37043076
$ws.Range("M122").Value = -37040626
$ws.Range("H132").Value = 1976.093
$ws.Range("I132").Value = 1958
$ws.Range("K132").Value = 5874
$ws.Range("M132").Value = -3344
$ws.Range("H135").Value = 72979.71000000001
$ws.Range("J135").Value = 72979.71000000001
$ws.Range("L135").Value = 72979.71000000001
$ws.Range("N135").Value = -83119.71000000001
$ws.Range("H136").Value = 4801.1875
$ws.Range("I136").Value = 4801.1875
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 14403.5625
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -11853.5625

# ===== BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1439.9
$ws.Range("I20").Value = 1550
$ws.Range("J20").Value = 999.5
$ws.Range("K20").Value = 1550
$ws.Range("L20").Value = 999.5
$ws.Range("M20").Value = -1303
$ws.Range("N20").Value = -1493.5
$ws.Range("H80").Value = 322.70587
$ws.Range("I80").Value = 305.36365
$ws.Range("J80").Value = 331
$ws.Range("K80").Value = 305.36365
$ws.Range("L80").Value = 331
$ws.Range("M80").Value = 692.63635
$ws.Range("N80").Value = -2327
$ws.Range("H83").Value = 322.70587
$ws.Range("I83").Value = 305.36365
$ws.Range("J83").Value = 331
$ws.Range("K83").Value = 1526.81825
$ws.Range("L83").Value = 1655
$ws.Range("M83").Value = 3465.18175
$ws.Range("N83").Value = -11639
$ws.Range("H105").Value = 2411.2896
$ws.Range("I105").Value = 2021.7241
$ws.Range("K105").Value = 2021.7241
$ws.Range("M105").Value = -274.7240999999999
$ws.Range("H107").Value = 2314.394
$ws.Range("I107").Value = 746.73914
$ws.Range("K107").Value = 746.73914
$ws.Range("M107").Value = 1173.26086
$ws.Range("H134").Value = 3012.353
$ws.Range("I134").Value = 3029.2856
$ws.Range("K134").Value = 9087.856800000001
$ws.Range("M134").Value = -6552.856800000001

# ===== CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2407.1538
$ws.Range("I58").Value = 2284.1738
$ws.Range("J58").Value = 3350
$ws.Range("K58").Value = 2284.1738
$ws.Range("L58").Value = 3350
$ws.Range("M58").Value = -2081.1738
$ws.Range("N58").Value = -3756
$ws.Range("H122").Value = 1021.4
$ws.Range("I122").Value = 988.5
$ws.Range("J122").Value = 1043.3334
$ws.Range("K122").Value = 2965.5
$ws.Range("L122").Value = 3130.0002
$ws.Range("M122").Value = -515.5
$ws.Range("N122").Value = -8030.0002
$ws.Range("H132").Value = 4547.4165
$ws.Range("I132").Value = 4547.4165
$ws.Range("K132").Value = 13642.2495
$ws.Range("M132").Value = -11112.2495
$ws.Range("H134").Value = 4644.5884
$ws.Range("I134").Value = 3871.5833
$ws.Range("K134").Value = 11614.7499
$ws.Range("M134").Value = -9079.749899999999
$ws.Range("H136").Value = 2407.1538
$ws.Range("I136").Value = 2284.1738
$ws.Range("J136").Value = 3350
$ws.Range("K136").Value = 6852.5214
$ws.Range("L136").Value = 10050
$ws.Range("M136").Value = -4302.5214
$ws.Range("N136").Value = -15150

# ===== CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 524.75
$ws.Range("I47").Value = 533
$ws.Range("K47").Value = 1599
$ws.Range("M47").Value = -1168

# ===== GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 9203.833000000001
$ws.Range("I97").Value = 768.75
$ws.Range("K97").Value = 768.75
$ws.Range("M97").Value = -272.75
$ws.Range("H122").Value = 12868.889
$ws.Range("I122").Value = 12708.826
$ws.Range("K122").Value = 38126.478
$ws.Range("M122").Value = -35676.478
$ws.Range("H126").Value = 2871.375
$ws.Range("I126").Value = 2839.5
$ws.Range("J126").Value = 3094.5
$ws.Range("K126").Value = 8518.5
$ws.Range("L126").Value = 9283.5
$ws.Range("M126").Value = -6048.5
$ws.Range("N126").Value = -14223.5
$ws.Range("H132").Value = 2501.7273
$ws.Range("I132").Value = 2307.7693
$ws.Range("K132").Value = 6923.3079
$ws.Range("M132").Value = -4393.3079

# ===== LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 917.88464
$ws.Range("I46").Value = 675.84
$ws.Range("J46").Value = 6969
$ws.Range("K46").Value = 675.84
$ws.Range("L46").Value = 6969
$ws.Range("M46").Value = -487.84
$ws.Range("N46").Value = -7345
$ws.Range("H61").Value = 2982.2903
$ws.Range("I61").Value = 1776.25
$ws.Range("K61").Value = 1776.25
$ws.Range("M61").Value = -1574.25
$ws.Range("H82").Value = 4999.857
$ws.Range("I82").Value = 2000
$ws.Range("K82").Value = 2000
$ws.Range("M82").Value = -1639
$ws.Range("H85").Value = 4999.857
$ws.Range("I85").Value = 2000
$ws.Range("K85").Value = 2000
$ws.Range("M85").Value = -752
$ws.Range("H113").Value = 2982.2903
$ws.Range("I113").Value = 1776.25
$ws.Range("K113").Value = 1776.25
$ws.Range("M113").Value = 393.75
$ws.Range("H132").Value = 3096.8386
$ws.Range("I132").Value = 3098.2856
$ws.Range("J132").Value = 3093.8
$ws.Range("K132").Value = 9294.856800000001
$ws.Range("L132").Value = 9281.400000000001
$ws.Range("M132").Value = -6764.856800000001
$ws.Range("N132").Value = -14341.4
$ws.Range("H136").Value = 10369.05
$ws.Range("I136").Value = 1155.9445
$ws.Range("K136").Value = 3467.8335
$ws.Range("M136").Value = -917.8335000000002

# ===== WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6587.1333
$ws.Range("I122").Value = 3633.6667
$ws.Range("J122").Value = 8556.111000000001
$ws.Range("K122").Value = 10901.0001
$ws.Range("L122").Value = 25668.333
$ws.Range("M122").Value = -8451.000100000001
$ws.Range("N122").Value = -30568.333
$ws.Range("H126").Value = 1864.2667
$ws.Range("I126").Value = 1361.8889
$ws.Range("K126").Value = 4085.6667
$ws.Range("M126").Value = -1615.6667
$ws.Range("H136").Value = 1485.4286
$ws.Range("I136").Value = 936
$ws.Range("K136").Value = 2808
$ws.Range("M136").Value = -258
